# --- Step 1: duplicate the pristine "Checklist" sheet to the end, to become "Question" ---
$wb = $excel.ActiveWorkbook
$checklist = $wb.Worksheets.Item("Checklist")
$timesheet = $wb.Worksheets.Item("Timesheet")
$ideas = $wb.Worksheets.Item("Ideas")

$checklist.Copy([System.Type]::Missing, $ideas)
$question = $wb.Worksheets.Item($wb.Worksheets.Count)
$question.Name = "Question"

# --- Step 2: edit Checklist sheet ---
# Header row (row 2): Item / Who? / Date accomplished -- copy header style from Ideas!C2
$ideas.Range("C2:E2").Copy()
$checklist.Range("B2:D2").PasteSpecial(-4122) # xlPasteFormats
$checklist.Range("B2").Value2 = "Item"
$checklist.Range("C2").Value2 = "Who?"
$checklist.Range("D2").Value2 = "Date accomplished"

# Row 3: Parse article files (HTML) / Eric Gan
$checklist.Rows.Item(3).RowHeight = 32.55
$checklist.Range("B3").Value2 = "Parse article files (HTML)"
$ideas.Range("C3").Copy()
$checklist.Range("C3").PasteSpecial(-4122) # xlPasteFormats
$checklist.Range("C3").Value2 = "Eric Gan"

# Row 4: Categorize question difficulty
$checklist.Rows.Item(4).RowHeight = 32.35
$checklist.Range("B4").Value2 = "Categorize question difficulty"

# Column D width
$checklist.Columns.Item(4).ColumnWidth = 14.43

# --- Step 3: edit Timesheet sheet ---
$ideas.Range("C3").Copy()
$timesheet.Range("C3").PasteSpecial(-4122) # xlPasteFormats
$timesheet.Range("C3").Value2 = "Downloaded NLTK"

$ideas.Range("C4").Copy()
$timesheet.Range("C4").PasteSpecial(-4122) # xlPasteFormats
$timesheet.Range("C4").Value2 = "Researched NLTK and its module for POS tagging"

# --- Step 4: edit Ideas sheet ---
$ideas.Rows.Item(4).RowHeight = 56.35
$ideas.Range("C4").Copy()
$ideas.Range("D4:E4").PasteSpecial(-4122) # xlPasteFormats
$ideas.Range("D4").Value2 = "Questions that might have a dependent clause appended to it."
$ideas.Range("E4").Value2 = "Ignore them for now"

# --- Step 5: finish the Question sheet ---
# Clear the bold header-label left over from the Checklist template (row 4, col B)
$question.Range("B5").Copy()
$question.Range("B4").PasteSpecial(-4122) # xlPasteFormats
$question.Range("B4").ClearContents()

$question.Columns.Item(3).ColumnWidth = 29.0

$question.Range("C3").Value2 = "Are we able to use/install nltk?"
$question.Range("C4").Value2 = "What is our data limit (in size)?"
$question.Range("C5").Value2 = "We need to install nltk_data onto linux machines."

Write-Host "Checklist B2:" $checklist.Range("B2").Value2
Write-Host "Checklist C2:" $checklist.Range("C2").Value2
Write-Host "Checklist D2:" $checklist.Range("D2").Value2
Write-Host "Checklist B3:" $checklist.Range("B3").Value2
Write-Host "Checklist C3:" $checklist.Range("C3").Value2
Write-Host "Checklist B4:" $checklist.Range("B4").Value2
Write-Host "Timesheet C3:" $timesheet.Range("C3").Value2
Write-Host "Timesheet C4:" $timesheet.Range("C4").Value2
Write-Host "Ideas D4:" $ideas.Range("D4").Value2
Write-Host "Ideas E4:" $ideas.Range("E4").Value2
Write-Host "Question B4:" $question.Range("B4").Value2
Write-Host "Question C3:" $question.Range("C3").Value2
Write-Host "Question C4:" $question.Range("C4").Value2
Write-Host "Question C5:" $question.Range("C5").Value2
